# The app used to hard-code the worksheet name ("Hoja1") when reading the
# uploaded attendance file. That requirement was dropped, so the sample
# workbook that exercises the new code path no longer needs a meaningful
# sheet name - rename it to something arbitrary to prove the lookup no
# longer depends on it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "dkjshfdsjhfdjks"

# Leave the view the way it was when the workbook was last saved: scrolled
# down so row 12 is at the top of the window, with the cursor resting one
# row below the last data row (A27).
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 1
$ws.Range("A27").Select()
